$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header style from G1 into H1 (bold/centered/bordered header look)
$ws.Range("G1").Copy($ws.Range("H1"))

# Set the new "Save" column header and data values
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
